{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change (see commit message / diff):\n//   \"Validate that the username is at least 4 characters.\"\n// becomes two bullet items:\n//   \"Validate that if the username is less than 4 characters, an error\n//    message appears.\"\n//   \"Validate that if the password is less than 8 characters, an error\n//    message appears.\"\n//\n// Both the exact run-splitting shown in the OOXML diff and the paragraph\n// (list-item) formatting are reproduced by replacing the paragraph with a\n// minimal Flat-OPC WordprocessingML fragment via Range.insertOoxml (the\n// Office.js equivalent of Word's InsertXML, which *replaces* the target\n// range's contents).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph to split/replace by its exact original text.\nconst targetText = \"Validate that the username is at least 4 characters.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find paragraph: \" + targetText);\n}\n\n// Shared run properties for every run in both paragraphs.\nconst RPR = '<w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>';\n// Shared paragraph properties (same bullet/list style as the surrounding items).\nconst PPR =\n  '<w:pPr><w:pStyle w:val=\"a3\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n  '<w:bidi w:val=\"0\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>';\n\nfunction run(text, preserve) {\n  const space = preserve ? ' xml:space=\"preserve\"' : \"\";\n  return \"<w:r>\" + RPR + \"<w:t\" + space + \">\" + text + \"</w:t></w:r>\";\n}\n\n// First (rewritten) paragraph: split across many runs exactly as in the diff.\nconst para1Runs = [\n  [\"Validate that \", true],\n  [\"if \", true],\n  [\"the user\", false],\n  [\"name is \", true],\n  [\"less than \", true],\n  [\"4 characters\", false],\n  [\", \", true],\n  [\"an e\", false],\n  [\"rror message appears\", false],\n  [\".\", false],\n];\n\n// Second (new) paragraph about the password, added right after the first.\nconst para2Runs = [\n  [\"Validate that if the password is less than 8 characters, \", true],\n  [\"an \", true],\n  [\"error message appears.\", false],\n];\n\nconst para1Xml = \"<w:p>\" + PPR + para1Runs.map(([t, p]) => run(t, p)).join(\"\") + \"</w:p>\";\nconst para2Xml = \"<w:p>\" + PPR + para2Runs.map(([t, p]) => run(t, p)).join(\"\") + \"</w:p>\";\n\nconst wordXml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  para1Xml +\n  para2Xml +\n  \"</w:body></w:document>\";\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  wordXml +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n// InsertXML REPLACES the addressed range's contents, so calling it on the\n// paragraph itself swaps that one paragraph for our two new paragraphs.\ntarget.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Target change (see commit message / diff):\n#   \"Validate that the username is at least 4 characters.\"\n# becomes two bullet items:\n#   \"Validate that if the username is less than 4 characters, an error\n#    message appears.\"\n#   \"Validate that if the password is less than 8 characters, an error\n#    message appears.\"\n#\n# We locate the paragraph by its exact original text, then use\n# Range.InsertXML to replace it with raw WordprocessingML for two\n# paragraphs, reproducing the exact run-splitting shown in the OOXML\n# diff as well as the shared bullet/list paragraph formatting.\n\n$d = $word.ActiveDocument\n\n$targetText = \"Validate that the username is at least 4 characters.\"\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    # Paragraph.Range.Text includes the trailing paragraph mark (CR); trim it\n    # (and any stray whitespace) before comparing against the plain text.\n    $paraText = $para.Range.Text.TrimEnd(\"`r\", \"`n\", \"`a\")\n    if ($paraText -eq $targetText) {\n        $target = $para\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find paragraph: $targetText\"\n}\n\n# Shared run properties for every run in both paragraphs.\n$RPR = '<w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>'\n# Shared paragraph properties (same bullet/list style as the surrounding items).\n$PPR = '<w:pPr><w:pStyle w:val=\"a3\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n       '<w:bidi w:val=\"0\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>'\n\nfunction New-Run([string]$text, [bool]$preserve) {\n    $space = \"\"\n    if ($preserve) { $space = ' xml:space=\"preserve\"' }\n    return \"<w:r>$RPR<w:t$space>$text</w:t></w:r>\"\n}\n\n# First (rewritten) paragraph: split across many runs exactly as in the diff.\n$para1Runs = @(\n    @(\"Validate that \", $true),\n    @(\"if \", $true),\n    @(\"the user\", $false),\n    @(\"name is \", $true),\n    @(\"less than \", $true),\n    @(\"4 characters\", $false),\n    @(\", \", $true),\n    @(\"an e\", $false),\n    @(\"rror message appears\", $false),\n    @(\".\", $false)\n)\n\n# Second (new) paragraph about the password, added right after the first.\n$para2Runs = @(\n    @(\"Validate that if the password is less than 8 characters, \", $true),\n    @(\"an \", $true),\n    @(\"error message appears.\", $false)\n)\n\n$para1Body = \"\"\nforeach ($run in $para1Runs) {\n    $para1Body += New-Run $run[0] $run[1]\n}\n\n$para2Body = \"\"\nforeach ($run in $para2Runs) {\n    $para2Body += New-Run $run[0] $run[1]\n}\n\n$wNs = ' xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$para1Xml = \"<w:p$wNs>$PPR$para1Body</w:p>\"\n$para2Xml = \"<w:p$wNs>$PPR$para2Body</w:p>\"\n\n# InsertXML REPLACES the addressed range's contents, so calling it on the\n# paragraph's range swaps that one paragraph for our two new paragraphs.\n$target.Range.InsertXML($para1Xml + $para2Xml)\n"}
